# Apply the update to the DaySale workbook:
#  - Row 10 (ANTI-COX II 15MG/3ML 6 AMP): ratio, sell price and txn-count columns change
#  - Row 26 (TELFAST): product renamed, price + sell price change
#  - Row 32 (سرنجات 3 سم): sell price + txn-count columns change
#  - Grand total (P36) updated
#  - Report generation timestamp (A37) updated

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-NumericLookingText {
    param(
        [string]$cellRef,
        [string]$newValue
    )
    # These cells carry a numeric display format (e.g. "0.00") even though the
    # sheet actually stores plain text in them. Assigning the string directly
    # would get auto-coerced into a real number, so briefly force a text
    # format, assign, then restore the original numeric format/style.
    $range = $ws.Range($cellRef)
    $origFormat = $range.NumberFormat
    $range.NumberFormat = "@"
    $range.Value = $newValue
    $range.NumberFormat = $origFormat
}

# --- Row 10: ANTI-COX II 15MG/3ML 6 AMP ---
$ws.Range("H10").Value = "1:0"
Set-NumericLookingText "P10" "117.0000"
$ws.Range("Q10").Value = "1:3"

# --- Row 26: TELFAST renamed + price updates ---
$ws.Range("C26").Value = "TELFAST 120MG 20 F.C. TAB"
$ws.Range("N26").Value = "99.00"
Set-NumericLookingText "P26" "49.5000"

# --- Row 32: سرنجات 3 سم ---
Set-NumericLookingText "P32" "18.0000"
$ws.Range("Q32").Value = "9:0"

# --- Grand total ---
$ws.Range("P36").Value = 1531.6700000000001

# --- Report timestamp ---
$ws.Range("A37").Value = "Monday, 29 September, 2025 12:04 PM"
